# Updated ITA model - add new "Electricity Trade Data" block and extra
# unnamed header columns to the historical_data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("historical_data")

# --- Extend header row (row 1) with "Unnamed: N" labels for columns I:Z ---
$unnamedCols = @("I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")
$unnamedIdx = 8
foreach ($col in $unnamedCols) {
    $ws.Range("$col`1").Value = "Unnamed: $unnamedIdx"
    $unnamedIdx++
}

# --- New block: Electricity Trade Data (TWh) ---
$ws.Range("A12").Value = "Electricity Trade Data (TWh)"

$ws.Range("A13").Value = "ISO"
$ws.Range("B13").Value = "attribute"

$years = @(2000,2001,2002,2003,2004,2005,2006,2007,2008,2009,2010,2011,2012,2013,2014,2015,2016,2017,2018,2019,2020,2021,2022,2023)
$yearCols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")
for ($i = 0; $i -lt $years.Length; $i++) {
    $ws.Range("$($yearCols[$i])13").Value = $years[$i]
}

$ws.Range("A14").Value = "ITA"
$ws.Range("B14").Value = "Export"
$ws.Range("C14").Value = 0

$ws.Range("A15").Value = "ITA"
$ws.Range("B15").Value = "Import"
$ws.Range("C15").Value = 44.35
